# Deep sea double count fix
# Updates "Status by Landings (Area)" and "Status by Landings (Tier)" sheets:
#  - MSF/O/Sustainable/Unsustainable (Mt) values get recalculated and a new
#    3-decimal custom number format (#,##0.000)
#  - MSF/O/Sustainable/Unsustainable (%) values get recalculated (format unchanged)

$wb = $excel.ActiveWorkbook

$newFmt = "#,##0.000"

# ---------------------------------------------------------------------------
# Sheet: "Status by Landings (Area)"
# ---------------------------------------------------------------------------
$wsArea = $wb.Worksheets.Item("Status by Landings (Area)")

# (Mt) values - rows 3-6, column C - also apply new number format
$wsArea.Range("C3").Value = 0.04298276364640884
$wsArea.Range("C3").NumberFormat = $newFmt

$wsArea.Range("C4").Value = 0.03504234798678382
$wsArea.Range("C4").NumberFormat = $newFmt

$wsArea.Range("C5").Value = 0.04298276364640884
$wsArea.Range("C5").NumberFormat = $newFmt

$wsArea.Range("C6").Value = 0.03504234798678382
$wsArea.Range("C6").NumberFormat = $newFmt

# (%) values - rows 8-11, column C - format stays the same
$wsArea.Range("C8").Value = 55.08837186735091
$wsArea.Range("C9").Value = 44.91162813264909
$wsArea.Range("C10").Value = 55.08837186735091
$wsArea.Range("C11").Value = 44.91162813264909

# ---------------------------------------------------------------------------
# Sheet: "Status by Landings (Tier)"
# ---------------------------------------------------------------------------
$wsTier = $wb.Worksheets.Item("Status by Landings (Tier)")

# Row 4 - "Deep Sea"
$wsTier.Range("C4").Value = 0.04298276364640884
$wsTier.Range("C4").NumberFormat = $newFmt

$wsTier.Range("D4").Value = 0.03504234798678382
$wsTier.Range("D4").NumberFormat = $newFmt

$wsTier.Range("E4").Value = 0.04298276364640884
$wsTier.Range("E4").NumberFormat = $newFmt

$wsTier.Range("F4").Value = 0.03504234798678382
$wsTier.Range("F4").NumberFormat = $newFmt

$wsTier.Range("H4").Value = 55.08837186735091
$wsTier.Range("I4").Value = 44.91162813264909
$wsTier.Range("J4").Value = 55.08837186735091
$wsTier.Range("K4").Value = 44.91162813264909

# Row 5 - "Global"
$wsTier.Range("C5").Value = 0.04298276364640884
$wsTier.Range("C5").NumberFormat = $newFmt

$wsTier.Range("D5").Value = 0.03504234798678382
$wsTier.Range("D5").NumberFormat = $newFmt

$wsTier.Range("E5").Value = 0.04298276364640884
$wsTier.Range("E5").NumberFormat = $newFmt

$wsTier.Range("F5").Value = 0.03504234798678382
$wsTier.Range("F5").NumberFormat = $newFmt

$wsTier.Range("H5").Value = 55.08837186735091
$wsTier.Range("I5").Value = 44.91162813264909
$wsTier.Range("J5").Value = 55.08837186735091
$wsTier.Range("K5").Value = 44.91162813264909
